{"js": "// Communication Management Plan Template - update the \"Group Meeting\" timing\n// cell from \"Roughly bi-weekly @ Sunday\" to \"Roughly twice a week@ Tuesday/Thursday\",\n// and normalize two other cells whose runs were split mid-word\n// (\"Dr. \" + \"Mohamed El-Darieby\" + \",\" and the \" \" + \"b\" inside \"bi-weekly\").\n\nconst body = context.document.body;\n\n// 1) \"Mentor Meeting\" row, Audience cell: merge \"Dr. \" / \"Mohamed El-Darieby\" / \",\"\n//    into a single run (no visible text change).\nconst drResults = body.search(\"Dr. Mohamed El-Darieby,\", { matchCase: true, matchWholeWord: false });\ndrResults.load(\"text\");\nawait context.sync();\nif (drResults.items.length > 0) {\n  drResults.items[0].insertText(\"Dr. Mohamed El-Darieby,\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"Mentor Meeting\" row, Timing cell: merge the \" \" and \"b\" runs inside\n//    \"Roughly bi-weekly @ Thursday\" into a single \" b\" run (no visible text change).\nconst thursdayResults = body.search(\"Roughly bi-weekly @ Thursday\", { matchCase: true, matchWholeWord: false });\nthursdayResults.load(\"text\");\nawait context.sync();\nif (thursdayResults.items.length > 0) {\n  const spaceB = thursdayResults.items[0].search(\" b\", { matchCase: true, matchWholeWord: false });\n  spaceB.load(\"text\");\n  await context.sync();\n  if (spaceB.items.length > 0) {\n    spaceB.items[0].insertText(\" b\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 3) \"Group Meeting\" row, Timing cell: update the schedule text itself.\nconst sundayResults = body.search(\"Roughly bi-weekly @ Sunday\", { matchCase: true, matchWholeWord: false });\nsundayResults.load(\"text\");\nawait context.sync();\nif (sundayResults.items.length > 0) {\n  sundayResults.items[0].insertText(\"Roughly twice a week@ Tuesday/Thursday\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Communication Management Plan Template - update the \"Group Meeting\" timing\n# cell from \"Roughly bi-weekly @ Sunday\" to \"Roughly twice a week@ Tuesday/Thursday\",\n# and normalize two other cells whose runs were split mid-word\n# (\"Dr. \" + \"Mohamed El-Darieby\" + \",\" and the \" \" + \"b\" inside \"bi-weekly\").\n\n$d = $word.ActiveDocument\n\n# Word constants (Find.Execute positional args):\n#   Wrap:    wdFindContinue = 1\n#   Replace: wdReplaceOne = 1, wdReplaceAll = 2\n#   NOTE: wdReplaceAll (2) replaces every match in the *entire story*, even when\n#   Execute is called on a narrowed Range - so wdReplaceOne (1) is used throughout,\n#   since every search string below is already unique in the document.\n\n# 1) \"Mentor Meeting\" row, Audience cell: merge \"Dr. \" / \"Mohamed El-Darieby\" / \",\"\n#    into a single run (no visible text change).\n$range1 = $d.Content\n$range1.Find.Execute(\"Dr. Mohamed El-Darieby,\", $false, $false, $false, $false, $false, $true, 1, $false, \"Dr. Mohamed El-Darieby,\", 1) | Out-Null\n\n# 2) \"Mentor Meeting\" row, Timing cell: merge the \" \" and \"b\" runs inside\n#    \"Roughly bi-weekly @ Thursday\" into a single \" b\" run (no visible text change).\n#    Scope the Find to just that cell's phrase first, then replace only \" b\" within it\n#    so the \"R\"/\"oughly\" split (different run formatting) is left untouched.\n#    MatchWholeWord must be $false here since \" b\" is a partial-word match.\n$outer = $d.Content\n$outer.Find.Execute(\"Roughly bi-weekly @ Thursday\") | Out-Null\n$inner = $d.Range($outer.Start, $outer.End)\n$inner.Find.Execute(\" b\", $false, $false, $false, $false, $false, $true, 1, $false, \" b\", 1) | Out-Null\n\n# 3) \"Group Meeting\" row, Timing cell: update the schedule text itself.\n$range3 = $d.Content\n$range3.Find.Execute(\"Roughly bi-weekly @ Sunday\", $false, $false, $false, $false, $false, $true, 1, $false, \"Roughly twice a week@ Tuesday/Thursday\", 1) | Out-Null\n"}
